# Applies the "seed.fruit.A" update: splits the F1 and F2 cross labels in
# column A into F1a/F1b and F2a/F2b respectively, and updates the sheet's
# active selection to the whole of column A.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4 & 5 were both labelled "F1" -- differentiate them.
$ws.Range("A4").Value = "F1a"
$ws.Range("A5").Value = "F1b"

# Rows 6 & 7 were both labelled "F2" -- differentiate them.
$ws.Range("A6").Value = "F2a"
$ws.Range("A7").Value = "F2b"

# Select the full column A, with the active cell at the top (A1).
$ws.Range("A1:A1048576").Select()
